$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing rows 66-123 down to 67-124.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the same record as the row that
# used to be at position 66 (now at 67), except for an updated Fecha (D)
# and Volumen (J).
$ws.Range("A66").Value = 4
$ws.Range("B66").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C66").Value = "Los Lagos"
$ws.Range("D66").Value = 44705
$ws.Range("E66").Value = 10
$ws.Range("F66").Value = 100112052
$ws.Range("G66").Value = "Albahaca"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 90
$ws.Range("K66").Value = 7000
$ws.Range("L66").Value = 7000
$ws.Range("M66").Value = 7000
$ws.Range("N66").Value = "$/docena de matas"
$ws.Range("O66").Value = "Región Metropolitana"
$ws.Range("P66").Value = 1167
$ws.Range("Q66").Value = 6
$ws.Range("R66").Value = "Hortaliza"

# Match the D column's date number format used by the rest of the column.
$ws.Range("D66").NumberFormat = $ws.Range("D67").NumberFormat
